# "Add Dmg Red B to buffs"
# - Renames three existing LR entries to DFLR_* (now eligible for the
#   "Dmg Red B" buff) and recalculates every Summon Rating that buff
#   affects.
# - Appends 7 new rows (ids 15-21) for cards newly tracked after the
#   buff was added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename cards that now carry the Dmg Red B buff ------------------
$ws.Cells.Item(3, 2).Value = "DFLR_PHY_Buu_Bois"
$ws.Cells.Item(4, 2).Value = "DFLR_STR_GT_Duo"
$ws.Cells.Item(5, 2).Value = "DFLR_TEQ_Fusion_Zamasu"

# --- Updated Summon Rating values (column D) for existing rows -------
# Values are entered with a leading apostrophe so they are stored as
# text (matching the workbook's existing convention of text-typed
# numeric-looking cells) rather than being auto-coerced to numbers.
$ws.Cells.Item(2, 4).Value = "'9.761091559169799"
$ws.Cells.Item(3, 4).Value = "'20.550341814468673"
$ws.Cells.Item(4, 4).Value = "'5.100269432844584"
$ws.Cells.Item(5, 4).Value = "'89.05775365126465"
$ws.Cells.Item(6, 4).Value = "'28.360934907238175"
$ws.Cells.Item(8, 4).Value = "'12.367395542094705"
$ws.Cells.Item(9, 4).Value = "'15.54309710321133"
$ws.Cells.Item(13, 4).Value = "'64.33532095603707"
$ws.Cells.Item(14, 4).Value = "'7.773118256406022"
$ws.Cells.Item(15, 4).Value = "'6.303014270008721"

# --- Append new rows 16-22 (IDs 15-21) --------------------------------
$newRows = @(
    @("15", "DF_INT_ToP_Androids",  "0", "28.0987740482617"),
    @("16", "BU_STR_Universe_2",    "0", "1.25"),
    @("17", "DF_AGL_Berserk_Kale",  "2", "7.571905174214976"),
    @("18", "DF_PHY_God_Goku",      "1", "14.912581234610627"),
    @("19", "DF_PHY_Kid_Goku",      "1", "29.12503072830121"),
    @("20", "DF_TEQ_Fusing_Kefla",  "0", "20.0"),
    @("21", "DF_STR_Costume_Videl", "1", "5.0")
)

$row = 16
foreach ($entry in $newRows) {
    # Copy the formatting (bold, border, centered alignment) used by the
    # rest of column A down onto the new ID cell.
    $ws.Range("A15").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = "'" + $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = "'" + $entry[2]
    $ws.Cells.Item($row, 4).Value = "'" + $entry[3]

    $row = $row + 1
}

$excel.CutCopyMode = 0
